{"js": "// Revisions on intro paragraph:\n// Remove the leading sentence \"I'm interested in a role building a\n// Developer Relations practice. \" from the SUMMARY section's first\n// paragraph, leaving the rest of the paragraph text unchanged.\n\nconst sentenceToRemove =\n  \"I\\u2019m interested in a role building a Developer Relations practice. \";\n\nconst body = context.document.body;\nconst results = body.search(sentenceToRemove, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Revisions on intro paragraph:\n# Remove the leading sentence \"I'm interested in a role building a\n# Developer Relations practice. \" from the SUMMARY section's first\n# paragraph, leaving the rest of the paragraph text unchanged.\n\n$d = $word.ActiveDocument\n\n# Build the sentence with a real right single quotation mark (U+2019),\n# matching the character used in the source document.\n$rsquo = [char]0x2019\n$sentenceToRemove = \"I\" + $rsquo + \"m interested in a role building a Developer Relations practice. \"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $sentenceToRemove\n$find.Replacement.Text = \"\"\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
